$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The three data rows (row 4 "per person", row 5 "removal, thousand tons" and
#    row 6 "average population") switch their numeric display format from two
#    decimal places ("0.00") to one decimal place ("0.0") for every year column
#    that already held a number.
$ws.Range("D4:O4").NumberFormat = "0.0"
$ws.Range("D5:O5").NumberFormat = "0.0"
$ws.Range("D6:O6").NumberFormat = "0.0"

# 2) A new column P is appended for year 2022, mirroring the formatting that
#    column O (2021) already has for every row of the table.
$ws.Range("O2").Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null

$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3").PasteSpecial(-4122) | Out-Null
$ws.Range("P3").Value = 2022

$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial(-4122) | Out-Null
$ws.Range("P4").Formula = "=P5/P6*1000"

$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial(-4122) | Out-Null
$ws.Range("P5").Value = 1339.6

$ws.Range("O6").Copy() | Out-Null
$ws.Range("P6").PasteSpecial(-4122) | Out-Null
$ws.Range("P6").Value = 6300.5

$ws.Application.CutCopyMode = $false

# 3) Move the active selection, matching where the editor's cursor ended up.
$ws.Range("S4").Select() | Out-Null
